$d = $word.ActiveDocument

# Locate the paragraph containing the distinctive "Django framework" text
# (the start of the block being removed) so we don't depend on hard-coded
# paragraph indices.
$findRng = $d.Content
$found = $findRng.Find.Execute("Django framework", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Django framework' paragraph."
}
$djangoPara = $findRng.Paragraphs.First
$djangoIdx = $djangoPara.Index

# The two paragraphs that must be removed entirely (the blank spacer
# paragraph and the "front end ... HTML3" paragraph) immediately follow
# the Django paragraph.
$spacerPara = $d.Paragraphs.Item($djangoIdx + 1)
$frontEndPara = $d.Paragraphs.Item($djangoIdx + 2)

# Delete the "front end" paragraph completely (text + its own paragraph
# mark) by removing the span from its start to the start of the
# following paragraph. This leaves the following paragraph untouched.
$afterFrontEnd = $d.Paragraphs.Item($djangoIdx + 3)
$d.Range($frontEndPara.Range.Start, $afterFrontEnd.Range.Start).Delete()

# Delete the blank spacer paragraph the same way.
$afterSpacer = $d.Paragraphs.Item($djangoIdx + 1)
$d.Range($spacerPara.Range.Start, $afterSpacer.Range.Start).Delete()

# Empty out the Django paragraph's run (remove its text and run
# formatting) while keeping the paragraph itself (and its paragraph
# mark/properties) intact.
$djangoParaNow = $d.Paragraphs.Item($djangoIdx)
$emptyRunXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
    "<w:pPr><w:pStyle w:val='Normal'/><w:jc w:val='both'/><w:rPr>" + `
    "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
    "<w:sz w:val='18'/><w:szCs w:val='18'/><w:lang w:val='en-IN'/></w:rPr></w:pPr>" + `
    "<w:r><w:rPr/></w:r></w:p>"
$djangoParaNow.Range.InsertXML($emptyRunXml)
